$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 38

$ws.Cells.Item($newRow, 1).Value = "2025-04-29 02:21:11"
$ws.Cells.Item($newRow, 2).Value = 93
